# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F) and "最低票价" (G) values to the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2894
$ws1.Range("G3").Value = 70
$ws1.Range("F5").Value = 6395
$ws1.Range("F6").Value = 2499
$ws1.Range("F10").Value = 2913
$ws1.Range("F11").Value = 349
$ws1.Range("F13").Value = 7235
$ws1.Range("F14").Value = 318
$ws1.Range("F15").Value = 44
$ws1.Range("F17").Value = 227
$ws1.Range("F20").Value = 8627
$ws1.Range("F28").Value = 94
$ws1.Range("F30").Value = 17
$ws1.Range("F31").Value = 41
$ws1.Range("F33").Value = 96
$ws1.Range("F34").Value = 2606
$ws1.Range("F40").Value = 707
$ws1.Range("F41").Value = 3726
$ws1.Range("F42").Value = 9
$ws1.Range("F43").Value = 186
$ws1.Range("F46").Value = 199
$ws1.Range("F47").Value = 36

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2894
$ws4.Range("G3").Value = 70
$ws4.Range("F6").Value = 6395
$ws4.Range("F7").Value = 2499
$ws4.Range("F12").Value = 2913
$ws4.Range("F13").Value = 349
$ws4.Range("F17").Value = 7235
$ws4.Range("F18").Value = 318
$ws4.Range("F20").Value = 227
$ws4.Range("F23").Value = 8627
$ws4.Range("F29").Value = 94
$ws4.Range("F31").Value = 41
$ws4.Range("F34").Value = 96
$ws4.Range("F35").Value = 2606
$ws4.Range("F40").Value = 707
$ws4.Range("F42").Value = 3726
$ws4.Range("F43").Value = 186
$ws4.Range("F47").Value = 199
$ws4.Range("F48").Value = 36
